# Remove the pantry ingredient rows that are no longer needed
# (garlic cloves, butter, onions, red onions), keeping olive oil,
# salt, vegetable oil and water with their original ingredient_id values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so row numbers for earlier deletions stay valid.
$ws.Rows(9).Delete()   # red onions
$ws.Rows(6).Delete()   # onions
$ws.Rows(4).Delete()   # butter
$ws.Rows(2).Delete()   # garlic cloves

# Match the saved selection state from the edited workbook.
$ws.Range("A4").Select()
